$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update A17 label: "V100 (16GB)" -> "8×V100 (16GB)"
$ws.Range("A17").Value = "8×V100 (16GB)"

# Update A18 label: "V100 (32GB)" -> "8×V100 (32GB)"
$ws.Range("A18").Value = "8×V100 (32GB)"

# Update E22 value: 104 -> 103
$ws.Range("E22").Value = 103

# Update A27 label: "Examples/s" -> "Examples/s on 8×V100 (16GB)"
$ws.Range("A27").Value = "Examples/s on 8×V100 (16GB)"

# Update A34 label: "Examples/s" -> "Examples/s on 8×V100 (16GB)"
$ws.Range("A34").Value = "Examples/s on 8×V100 (16GB)"

# Update sheet view: topLeftCell B33->G20, zoom 130->140, selection D40->E22
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.Zoom = 140
$ws.Range("E22").Select()
